$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305475234985352
$ws.Range("B1").Value = 3.813564777374268
$ws.Range("C1").Value = 3.86201548576355
$ws.Range("D1").Value = 2.997583627700806
$ws.Range("E1").Value = 1.045259952545166
